# Auto-generated script applying numeric updates to Cactuar_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2106270
$ws.Range("I40").Value = 18499.834
$ws.Range("J40").Value = 3358932.2
$ws.Range("K40").Value = 18499.834
$ws.Range("L40").Value = 3358932.2
$ws.Range("M40").Value = -18324.834
$ws.Range("N40").Value = -3359282.2
$ws.Range("H43").Value = 3080498.8
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5138
$ws.Range("H51").Value = 6753.478
$ws.Range("J51").Value = 7430
$ws.Range("L51").Value = 7430
$ws.Range("N51").Value = -8398
$ws.Range("H70").Value = 6206.4165
$ws.Range("I70").Value = 3812.8333
$ws.Range("J70").Value = 8600
$ws.Range("K70").Value = 11438.4999
$ws.Range("L70").Value = 25800
$ws.Range("M70").Value = -11168.4999
$ws.Range("N70").Value = -26340
$ws.Range("H73").Value = 6206.4165
$ws.Range("I73").Value = 3812.8333
$ws.Range("J73").Value = 8600
$ws.Range("K73").Value = 11438.4999
$ws.Range("L73").Value = 25800
$ws.Range("M73").Value = -10502.4999
$ws.Range("N73").Value = -27672
$ws.Range("H112").Value = 2994.327
$ws.Range("J112").Value = 3071.898
$ws.Range("L112").Value = 9215.694
$ws.Range("N112").Value = -11431.694
$ws.Range("H132").Value = 6419.772
$ws.Range("I132").Value = 2337.2063
$ws.Range("K132").Value = 7011.618899999999
$ws.Range("M132").Value = -4481.618899999999
$ws.Range("H137").Value = 15154137
$ws.Range("J137").Value = 37039836
$ws.Range("L137").Value = 111119508
$ws.Range("N137").Value = -111124608
$ws.Range("H138").Value = 5629.8535
$ws.Range("I138").Value = 2494.7273
$ws.Range("J138").Value = 6931.2266
$ws.Range("K138").Value = 7484.1819
$ws.Range("L138").Value = 20793.6798
$ws.Range("M138").Value = -2344.1819
$ws.Range("N138").Value = -31073.6798
$ws.Range("H141").Value = 3670.1738
$ws.Range("I141").Value = 3918.0952
$ws.Range("K141").Value = 11754.2856
$ws.Range("M141").Value = -6574.285600000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14509.656
$ws.Range("I32").Value = 13993.526
$ws.Range("K32").Value = 13993.526
$ws.Range("M32").Value = -13706.526
$ws.Range("H61").Value = 4514.1377
$ws.Range("I61").Value = 3784.8948
$ws.Range("J61").Value = 5899.7
$ws.Range("K61").Value = 3784.8948
$ws.Range("L61").Value = 5899.7
$ws.Range("M61").Value = -3572.8948
$ws.Range("N61").Value = -6323.7
$ws.Range("H74").Value = 8622473
$ws.Range("I74").Value = 10418181
$ws.Range("K74").Value = 10418181
$ws.Range("M74").Value = -10417307
$ws.Range("H77").Value = 8622473
$ws.Range("I77").Value = 10418181
$ws.Range("K77").Value = 52090905
$ws.Range("M77").Value = -52086537
$ws.Range("H132").Value = 27127.74
$ws.Range("I132").Value = 35485.234
$ws.Range("K132").Value = 106455.702
$ws.Range("M132").Value = -103925.702
$ws.Range("H136").Value = 4514.1377
$ws.Range("I136").Value = 3784.8948
$ws.Range("J136").Value = 5899.7
$ws.Range("K136").Value = 11354.6844
$ws.Range("L136").Value = 17699.1
$ws.Range("M136").Value = -8804.6844
$ws.Range("N136").Value = -22799.1
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 78527280
$ws.Range("J99").Value = 4100
$ws.Range("L99").Value = 4100
$ws.Range("N99").Value = -7096
$ws.Range("H107").Value = 1929
$ws.Range("I107").Value = 1934.8
$ws.Range("K107").Value = 1934.8
$ws.Range("M107").Value = -14.79999999999995
$ws.Range("H134").Value = 2969.7144
$ws.Range("I134").Value = 2969.7144
$ws.Range("K134").Value = 8909.143199999999
$ws.Range("M134").Value = -6374.143199999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 70619.625
$ws.Range("J74").Value = 69381.71000000001
$ws.Range("L74").Value = 69381.71000000001
$ws.Range("N74").Value = -71129.71000000001
$ws.Range("H77").Value = 70619.625
$ws.Range("J77").Value = 69381.71000000001
$ws.Range("L77").Value = 208145.13
$ws.Range("N77").Value = -216881.13
$ws.Range("H94").Value = 2171.762
$ws.Range("I94").Value = 905.6
$ws.Range("J94").Value = 2567.4375
$ws.Range("K94").Value = 905.6
$ws.Range("L94").Value = 2567.4375
$ws.Range("M94").Value = -454.6
$ws.Range("N94").Value = -3469.4375
$ws.Range("H99").Value = 8233.8125
$ws.Range("I99").Value = 7905.125
$ws.Range("K99").Value = 7905.125
$ws.Range("M99").Value = -6407.125
$ws.Range("H105").Value = 1625132.1
$ws.Range("I105").Value = 2526724.2
$ws.Range("K105").Value = 2526724.2
$ws.Range("M105").Value = -2524977.2
$ws.Range("H126").Value = 8233.8125
$ws.Range("I126").Value = 7905.125
$ws.Range("K126").Value = 23715.375
$ws.Range("M126").Value = -21245.375
$ws.Range("H132").Value = 35101708
$ws.Range("I132").Value = 41676410
$ws.Range("J132").Value = 36632.668
$ws.Range("K132").Value = 125029230
$ws.Range("L132").Value = 109898.004
$ws.Range("M132").Value = -125026700
$ws.Range("N132").Value = -114958.004
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 83.375
$ws.Range("I8").Value = 83.375
$ws.Range("K8").Value = 250.125
$ws.Range("M8").Value = -111.125
$ws.Range("H51").Value = 4999.5
$ws.Range("J51").Value = 4999.5
$ws.Range("L51").Value = 14998.5
$ws.Range("N51").Value = -15918.5
$ws.Range("H132").Value = 6707.4546
$ws.Range("I132").Value = 1160.1666
$ws.Range("K132").Value = 10441.4994
$ws.Range("M132").Value = -7911.499400000001
$ws.Range("H133").Value = 22098.834
$ws.Range("I133").Value = 9197.666999999999
$ws.Range("J133").Value = 35000
$ws.Range("K133").Value = 27593.001
$ws.Range("L133").Value = 105000
$ws.Range("M133").Value = -22533.001
$ws.Range("N133").Value = -115120
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50302
$ws.Range("H46").Value = 39000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640
$ws.Range("H70").Value = 1981812.4
$ws.Range("I70").Value = 2397299.2
$ws.Range("J70").Value = 8250
$ws.Range("K70").Value = 2397299.2
$ws.Range("L70").Value = 8250
$ws.Range("M70").Value = -2397029.2
$ws.Range("N70").Value = -8790
$ws.Range("H73").Value = 1981812.4
$ws.Range("I73").Value = 2397299.2
$ws.Range("J73").Value = 8250
$ws.Range("K73").Value = 2397299.2
$ws.Range("L73").Value = 8250
$ws.Range("M73").Value = -2396363.2
$ws.Range("N73").Value = -10122
$ws.Range("H80").Value = 1863833
$ws.Range("I80").Value = 2774916.2
$ws.Range("K80").Value = 2774916.2
$ws.Range("M80").Value = -2773918.2
$ws.Range("H83").Value = 1863833
$ws.Range("I83").Value = 2774916.2
$ws.Range("K83").Value = 13874581
$ws.Range("M83").Value = -13869589
$ws.Range("H102").Value = 26326756
$ws.Range("J102").Value = 8331.299999999999
$ws.Range("L102").Value = 8331.299999999999
$ws.Range("N102").Value = -11575.3
$ws.Range("H122").Value = 339888.7
$ws.Range("I122").Value = 850461.3
$ws.Range("K122").Value = 2551383.9
$ws.Range("M122").Value = -2548933.9
$ws.Range("H132").Value = 4797
$ws.Range("I132").Value = 4507.9443
$ws.Range("K132").Value = 13523.8329
$ws.Range("M132").Value = -10993.8329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2094.1924
$ws.Range("I16").Value = 1215.0869
$ws.Range("J16").Value = 8834
$ws.Range("K16").Value = 1215.0869
$ws.Range("L16").Value = 8834
$ws.Range("M16").Value = -1045.0869
$ws.Range("N16").Value = -9174
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 11995
$ws.Range("J19").Value = 11995
$ws.Range("L19").Value = 11995
$ws.Range("N19").Value = -12343
$ws.Range("H126").Value = 7082.909
$ws.Range("I126").Value = 6500
$ws.Range("J126").Value = 7301.5
$ws.Range("K126").Value = 19500
$ws.Range("L126").Value = 21904.5
$ws.Range("M126").Value = -17030
$ws.Range("N126").Value = -26844.5
$ws.Range("H136").Value = 8358.459000000001
$ws.Range("I136").Value = 3640.6216
$ws.Range("J136").Value = 11220.099
$ws.Range("K136").Value = 10921.8648
$ws.Range("L136").Value = 33660.297
$ws.Range("M136").Value = -8371.864799999999
$ws.Range("N136").Value = -38760.297

Write-Host "Applied 218 cell updates"